# Update the cryptocurrency price/volume table with freshly scraped values.
# Column D ("Price") holds numeric-looking text (e.g. "246.32", "0.999") that
# must stay plain text, exactly like the rest of the sheet -- Excel's COM
# Range.Value setter auto-coerces bare numeric strings into real numbers, so
# we temporarily force a Text number format on the whole Price column before
# writing, then restore the default ("Normal") style afterwards so the
# on-disk cell style indices are not perturbed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "35.788.85"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "1.891.03"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "246.32"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "0.693"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "43.09"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").Value = "56.41"
$ws.Range("E9").Value = "  +8.65%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").Value = "0.0754"
$ws.Range("E11").Value = "  +2.37%  "
$ws.Range("D12").Value = "0.0986"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "14.97"
$ws.Range("E13").Value = "  +14.23%  "
$ws.Range("D14").Value = "0.792"
$ws.Range("E14").Value = "  +7.76%  "
$ws.Range("D15").Value = "2.166.65"
$ws.Range("E15").Value = "  +0.73%  "
$ws.Range("D16").Value = "5.06"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "1.895.46"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "35.764.16"
$ws.Range("D19").Value = "73.63"
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "247.18"
$ws.Range("E21").Value = "  +0.68%  "
$ws.Range("D22").Value = "13.09"
$ws.Range("E22").Value = "  +2.30%  "
$ws.Range("D23").Value = "5.17"
$ws.Range("E23").Value = "  +4.80%  "
$ws.Range("E24").Value = "  +4.44%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -1.32%  "
$ws.Range("D27").Value = "166.11"
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "8.66"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").Value = "18.44"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "4.44"
$ws.Range("E31").Value = "  +4.69%  "
$ws.Range("D32").Value = "0.0611"
$ws.Range("E32").Value = "  +5.56%  "
$ws.Range("D33").Value = "4.27"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").Value = "1.85"
$ws.Range("E34").Value = "  +18.84%  "
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("E36").Value = "  -14.29%  "
$ws.Range("D37").Value = "0.856"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").Value = "0.0774"
$ws.Range("E38").Value = "  +12.38%  "
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("E40").Value = "  +6.99%  "
$ws.Range("D41").Value = "99.05"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("D42").Value = "16.98"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("D44").Value = "14.11"
$ws.Range("E44").Value = "  +16.48%  "
$ws.Range("D45").Value = "1.313.14"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("D50").Value = "6.32"
$ws.Range("E50").Value = "  +0.92%  "
$ws.Range("D51").Value = "42.62"
$ws.Range("E51").Value = "  -1.25%  "

# Restore the default cell style on the Price column so we don't leave
# behind a lingering "@" (Text) number-format style that wasn't present in
# the original file.
$ws.Range("D2:D51").Style = "Normal"
